# LOM3058.xlsx edit: remove the stand-alone "Docentes responsáveis" value row
# (old row 13, which only held "1033242 - Fábio Herbst Florenzano" in B/C),
# shifting all following rows up by one, and then correct a handful of
# mis-matched B/C text values further down the sheet to their new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the old row 13 entirely (shifts rows 14-25 up to 13-24).
$ws.Rows.Item(13).Delete()

# 2) After the shift, patch the B/C cell text that changed content
#    (row numbers below are the NEW row numbers, post-delete).

# Row 13 "Programa resumido:" -> "Semestral"
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 "Programa:" -> "01/01/2020"
# (Use a formula-then-paste-values round trip through a scratch cell so the
#  COM layer keeps this a literal text string instead of auto-converting it
#  to a date serial number.)
$ws.Range("Z1").Formula = "=""01/01/2020"""
$ws.Range("Z1").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("Z1").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

# Row 18 "Método:" -> "1033242 - Fábio Herbst Florenzano"
$ws.Range("B18:C18").Value = "1033242 - Fábio Herbst Florenzano"

# Row 19 "Critério:" -> "Provas escritas envolvendo o conteúdo teórico ministrado em sala de aula."
$ws.Range("B19:C19").Value = "Provas escritas envolvendo o conteúdo teórico ministrado em sala de aula."

# Row 20 "Norma de recuperação:" -> "Duas avaliações, sendo que a nota final corresponde..."
$ws.Range("B20:C20").Value = "Duas avaliações, sendo que a nota final corresponde à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados. Alunos com notas finais situadas no intervalo de 3 a 5 serão encaminhados à recuperação."

# Row 21 "Bibliografia:" -> "O aluno será submetido a um programa de estudos..."
$ws.Range("B21:C21").Value = "O aluno será submetido a um programa de estudos destinado a rever o conteúdo abordado na disciplina. Ao final deste período será aplicada uma nova avaliação. A nota final do aluno será a média aritmética desta avaliação com a nota anteriormente obtida, estando aprovados os alunos que tiverem nota final igual ou superior a 5."
